$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in I1/J1, matching the formatting of the existing H1 header
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill the new data columns for rows 2-33:
#   I column = 1 (constant)
#   J column = same value as column H on that row
for ($r = 2; $r -le 33; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value2
}
